$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A (Beteckning) values.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 499 }

# Column C ("Förändrad") holds a date serial value that was bumped from
# 45188 (2023-09-19) to 45189 (2023-09-20) for every data row (rows 2..last).
$oldValue = 45188
$newValue = 45189

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
